$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201 (pushes former rows 201-226 down to 202-227)
$ws.Rows.Item(201).Insert()

# Fill the new row 201 with the new weekly record
$ws.Cells.Item(201, 1).Value = 6
$ws.Cells.Item(201, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(201, 3).Value = "Metropolitana"
$ws.Cells.Item(201, 4).Value = 44748
$ws.Cells.Item(201, 5).Value = 13
$ws.Cells.Item(201, 6).Value = 100112022
$ws.Cells.Item(201, 7).Value = "Arveja Verde"
$ws.Cells.Item(201, 8).Value = "Perfection"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 200
$ws.Cells.Item(201, 11).Value = 38000
$ws.Cells.Item(201, 12).Value = 40000
$ws.Cells.Item(201, 13).Value = 38800
$ws.Cells.Item(201, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(201, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(201, 16).Value = 1552
$ws.Cells.Item(201, 17).Value = 25
$ws.Cells.Item(201, 18).Value = "Hortaliza"

# Insert a second new row at 218 (after the first insert, the former row 217
# - now at 217 - is immediately above this new gap; former rows 217-226,
# now at 218-227, get pushed down to 219-228)
$ws.Rows.Item(218).Insert()

# Fill the new row 218 with the second new weekly record
$ws.Cells.Item(218, 1).Value = 6
$ws.Cells.Item(218, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(218, 3).Value = "Metropolitana"
$ws.Cells.Item(218, 4).Value = 44747
$ws.Cells.Item(218, 5).Value = 13
$ws.Cells.Item(218, 6).Value = 100112022
$ws.Cells.Item(218, 7).Value = "Arveja Verde"
$ws.Cells.Item(218, 8).Value = "Perfection"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 400
$ws.Cells.Item(218, 11).Value = 38000
$ws.Cells.Item(218, 12).Value = 40000
$ws.Cells.Item(218, 13).Value = 38850
$ws.Cells.Item(218, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(218, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(218, 16).Value = 1554
$ws.Cells.Item(218, 17).Value = 25
$ws.Cells.Item(218, 18).Value = "Hortaliza"
